# Updates cryptos list data per commit: Updated cryptos list on Mon Aug  7 21:25:25 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'29.168.59"
$ws.Range("E2").Value = "  +0.26%  "

# Row 3
$ws.Range("D3").Value = "'1.825.18"
$ws.Range("E3").Value = "  -0.31%  "

# Row 4
$ws.Range("D4").Value = "'0.9985"
$ws.Range("E4").Value = "  +0.24%  "

# Row 5
$ws.Range("D5").Value = "'241.57"
$ws.Range("E5").Value = "  -0.52%  "

# Row 6
$ws.Range("D6").Value = "'0.6176"
$ws.Range("E6").Value = "  -1.80%  "

# Row 7
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.27%  "

# Row 8
$ws.Range("D8").Value = "'0.07337"
$ws.Range("E8").Value = "  -2.94%  "

# Row 9
$ws.Range("D9").Value = "'0.2890"
$ws.Range("E9").Value = "  -1.31%  "

# Row 10
$ws.Range("D10").Value = "'22.96"
$ws.Range("E10").Value = "  -1.28%  "

# Row 11
$ws.Range("E11").Value = "  +0.03%  "

# Row 12
$ws.Range("D12").Value = "'1.823.00"
$ws.Range("E12").Value = "  -0.41%  "

# Row 13
$ws.Range("E13").Value = "  -1.36%  "

# Row 14
$ws.Range("D14").Value = "'0.6609"
$ws.Range("E14").Value = "  -1.17%  "

# Row 15
$ws.Range("D15").Value = "'82.04"

# Row 16
$ws.Range("D16").Value = "'0.000008910"
$ws.Range("E16").Value = "  -5.29%  "

# Row 17
$ws.Range("D17").Value = "'5.834"
$ws.Range("E17").Value = "  -2.80%  "

# Row 18
$ws.Range("D18").Value = "'29.139.00"
$ws.Range("E18").Value = "  +0.17%  "

# Row 19
$ws.Range("D19").Value = "'2.065.90"
$ws.Range("E19").Value = "  -0.63%  "

# Row 20
$ws.Range("D20").Value = "'237.76"
$ws.Range("E20").Value = "  +6.45%  "

# Row 21
$ws.Range("E21").Value = "  -1.44%  "

# Row 22
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  +0.13%  "

# Row 23
$ws.Range("D23").Value = "'7.170"
$ws.Range("E23").Value = "  +0.39%  "

# Row 24
$ws.Range("D24").Value = "'0.9996"
$ws.Range("E24").Value = "  +0.21%  "

# Row 25
$ws.Range("D25").Value = "'158.14"
$ws.Range("E25").Value = "  -0.94%  "

# Row 26
$ws.Range("E26").Value = "  +1.29%  "

# Row 27
$ws.Range("D27").Value = "'8.433"
$ws.Range("E27").Value = "  -0.89%  "

# Row 28
$ws.Range("D28").Value = "'17.63"
$ws.Range("E28").Value = "  -1.42%  "

# Row 29
$ws.Range("D29").Value = "'1.483"
$ws.Range("E29").Value = "  -0.66%  "

# Row 30
$ws.Range("D30").Value = "'0.05554"
$ws.Range("E30").Value = "  -3.13%  "

# Row 31
$ws.Range("D31").Value = "'4.091"
$ws.Range("E31").Value = "  -0.28%  "

# Row 32
$ws.Range("E32").Value = "  -1.55%  "

# Row 33
$ws.Range("D33").Value = "'1.204"
$ws.Range("E33").Value = "  +0.30%  "

# Row 34
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "'1.821"
$ws.Range("E34").Value = "  -0.71%  "

# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.7350"
$ws.Range("E35").Value = "  -1.01%  "

# Row 36
$ws.Range("E36").Value = "  -0.63%  "

# Row 37
$ws.Range("D37").Value = "'2.614"
$ws.Range("E37").Value = "  -1.77%  "

# Row 38
$ws.Range("D38").Value = "'2.834"
$ws.Range("E38").Value = "  +2.81%  "

# Row 39
$ws.Range("D39").Value = "'1.212.96"
$ws.Range("E39").Value = "  -1.34%  "

# Row 40
$ws.Range("D40").Value = "'0.01757"
$ws.Range("E40").Value = "  -1.27%  "

# Row 41
$ws.Range("E41").Value = "  -2.75%  "

# Row 42
$ws.Range("D42").Value = "'0.9195"
$ws.Range("E42").Value = "  +3.47%  "

# Row 43
$ws.Range("D43").Value = "'1.001"
$ws.Range("E43").Value = "  +0.27%  "

# Row 44
$ws.Range("D44").Value = "'0.00000000130"
$ws.Range("E44").Value = "  +6.53%  "

# Row 45
$ws.Range("E45").Value = "  -0.68%  "

# Row 46
$ws.Range("D46").Value = "'1.971.64"
$ws.Range("E46").Value = "  -0.41%  "

# Row 47
$ws.Range("D47").Value = "'64.63"
$ws.Range("E47").Value = "  -1.87%  "

# Row 48
$ws.Range("D48").Value = "'0.5080"
$ws.Range("E48").Value = "  -0.07%  "

# Row 49
$ws.Range("D49").Value = "'0.4005"
$ws.Range("E49").Value = "  -1.61%  "

# Row 50
$ws.Range("D50").Value = "'9.059"
$ws.Range("E50").Value = "  +0.74%  "

# Row 51
$ws.Range("D51").Value = "'0.05765"
$ws.Range("E51").Value = "  -0.89%  "
